# Rename the inline picture "name" metadata (wp:docPr / Name) for the
# three logo images embedded in this document's headers/footers:
#   - footer1 Pearson logo : image2.png -> image1.png
#   - header2 BTec logo    : image1.jpg -> image2.jpg
#   - footer2 Pearson logo : image2.png -> image1.png
#
# The images are located in Section 1's headers/footers, so we walk
# every Header/Footer slot and match on the picture's AlternativeText
# (the "descr" attribute) plus its current Name, rather than hard-coding
# positional indices, so the script is resilient to slot ordering.

$d = $word.ActiveDocument

function Rename-LogoPicture($headerFooter, $expectedAlt, $oldName, $newName) {
    if (-not $headerFooter.Exists) {
        return
    }
    $shapes = $headerFooter.Range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $pic = $shapes.Item($i)
        if ($pic.AlternativeText -eq $expectedAlt -and $pic.Name -eq $oldName) {
            $pic.Name = $newName
        }
    }
}

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers.Item($i)
        Rename-LogoPicture $hdr "BTec_Logo-Orange" "image1.jpg" "image2.jpg"

        $ftr = $sec.Footers.Item($i)
        Rename-LogoPicture $ftr "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" "image2.png" "image1.png"
    }
}
